$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the "Total of states" totals row
$ws.Range("A53:B53").ClearContents()

# 2. Header cell A1 loses its bold emphasis and its explicit left alignment
#    (reverts to the default/general alignment), while keeping its bottom border.
$ws.Range("A1").Font.Bold = $false
$ws.Range("A1").HorizontalAlignment = 1  # xlGeneral

# 3. Set column C width (closest this engine's 1/6-char quantization allows to the
#    target raw width of 12.578125; 11.6-11.74 all round to a stored width of 12.5)
$ws.Columns("C").ColumnWidth = 11.67

# 4. Move selection to I10 (also resets any scrolled viewport)
$ws.Range("I10").Select() | Out-Null
